$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Max "Hoots" Hooton - Unity developer
$ws.Range("D14").Value = "Max ""Hoots"" Hooton"
$ws.Range("E14").Value = "Unity developer"

# Row 15: Zech Birkel - QA
$ws.Range("D15").Value = "Zech Birkel"
$ws.Range("E15").Value = "QA"

# Row 16: Shane "Vhespir" - Music
$ws.Range("D16").Value = "Shane ""Vhespir"""
$ws.Range("E16").Value = "Music"

# Match the author's final selection recorded in the sheet view
$ws.Range("E15").Select() | Out-Null
